# daily auto push: 2026-01-06 02:26 UTC
# Insert a new data point (2026/01/06, 9:00, rank 158) at row 566 of the
# log sheet, shifting the existing rows (566-607, the 2026/12/29 ..
# 2027/01/05 entries) down by one row to 567-608.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every row from 566 onward down by one to make room for the new entry.
$ws.Rows.Item(566).Insert()

# Force column A to Text so the "yyyy/mm/dd" string is not auto-converted
# into a date serial number by Excel's input parsing (matches how the
# rest of column A is stored, as plain text).
$ws.Range("A566").NumberFormat = "@"

$ws.Range("A566").Value = "2026/01/06"
$ws.Range("B566").Value = "火"
$ws.Range("C566").Value = 9
$ws.Range("D566").Value = 158
